$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '43.369.72'
$ws.Range("E2").Value = '  -0.57%  '
$ws.Range("D3").Value = '2.281.42'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '112.63'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '265.58'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.98%  '
$ws.Range("E7").Value = '  +4.15%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '46.89'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.75%  '
$ws.Range("E11").Value = '  -1.09%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.32'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +2.31%  '
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.33'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.25%  '
$ws.Range("D15").Value = '2.622.75'
$ws.Range("E15").Value = '  -0.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.864'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +1.89%  '
$ws.Range("D17").Value = '2.278.35'
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("D18").Value = '43.344.51'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.75'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.23'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.99'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.00%  '
$ws.Range("E24").Value = '  +4.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.43'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.93%  '
$ws.Range("E26").Value = '  +1.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.44'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.07'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.35'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.17%  '
$ws.Range("E30").Value = '  -1.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '173.36'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.64'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +0.30%  '
$ws.Range("E33").Value = '  -3.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.65'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.30%  '
$ws.Range("E35").Value = '  +3.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0380'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.21%  '
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.90'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +2.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.105'
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.60'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +8.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.30'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.71'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.238'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.08'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("E46").Value = '  -1.31%  '
$ws.Range("E47").Value = '  +3.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.58'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0998'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.63%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '100.33'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -2.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.614'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +12.64%  '
